# Add a new "2022" year column (S) to the Tourism-share-of-GDP table and
# refresh the last three years' data values, matching the authored diff:
#   - S4 = 2022                (header, same formatting as the other year cells)
#   - S5 = 3.4                 (data value, same formatting as the other data cells)
#   - P5/Q5/R5 values updated
#   - selection moves to T4 (one cell to the right of the newly added column)
#   - dimension/row spans are recalculated by Excel as part of normal editing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell S4 (copy formatting from the neighboring year header R4) ---
$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("S4").Value = 2022

# --- New data cell S5 (copy formatting from the neighboring data cell R5) ---
$ws.Range("R5").Copy() | Out-Null
$ws.Range("S5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("S5").Value = 3.4

$excel.CutCopyMode = 0

# --- Refresh the trailing data values for the existing years ---
$ws.Range("P5").Value = 4.4000000000000004
$ws.Range("Q5").Value = 2.9
$ws.Range("R5").Value = 3.2

# --- Move the active selection to T4, as recorded in the saved view state ---
$ws.Range("T4").Select() | Out-Null
